# Applies row permutation fix for CryCompanywiseStockReport: for several groups of
# rows sharing the same product (rows identified by the A column), the per-row
# B (Sauda/voucher no.), E (rate), F (qty) and G (value) figures were associated
# with the wrong row. This restores the correct values for each affected cell.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 142
$ws.Range("B142").Value = 48654
$ws.Range("E142").Value = 38.26
$ws.Range("F142").Value = -1
$ws.Range("G142").Value = -32.02

# Row 143
$ws.Range("B143").Value = 63902
$ws.Range("E143").Value = 34.04
$ws.Range("F143").Value = 2
$ws.Range("G143").Value = 64.04000000000001

# Row 154
$ws.Range("B154").Value = 53925
$ws.Range("E154").Value = 79.37
$ws.Range("F154").Value = 1
$ws.Range("G154").Value = 66.44

# Row 155
$ws.Range("B155").Value = 64350
$ws.Range("E155").Value = 70.63
$ws.Range("F155").Value = 101
$ws.Range("G155").Value = 6710.44

# Row 156
$ws.Range("B156").Value = 57756
$ws.Range("F156").Value = -100
$ws.Range("G156").Value = -6644

# Row 176
$ws.Range("B176").Value = 64329
$ws.Range("E176").Value = 128.32
$ws.Range("F176").Value = 6
$ws.Range("G176").Value = 724.14

# Row 177
$ws.Range("B177").Value = 57552
$ws.Range("E177").Value = 136.86
$ws.Range("F177").Value = -5
$ws.Range("G177").Value = -603.45

# Row 271
$ws.Range("B271").Value = 64973
$ws.Range("E271").Value = 35.4
$ws.Range("F271").Value = 150
$ws.Range("G271").Value = 4995

# Row 272
$ws.Range("B272").Value = 48706
$ws.Range("E272").Value = 39.8
$ws.Range("F272").Value = -144
$ws.Range("G272").Value = -4795.2

# Row 305
$ws.Range("B305").Value = 57854
$ws.Range("F305").Value = 2
$ws.Range("G305").Value = 611.6799999999999

# Row 306
$ws.Range("B306").Value = 62997
$ws.Range("F306").Value = 72
$ws.Range("G306").Value = 22020.48

# Row 308
$ws.Range("B308").Value = 63565
$ws.Range("E308").Value = 109.19
$ws.Range("F308").Value = 60
$ws.Range("G308").Value = 6162.6

# Row 309
$ws.Range("B309").Value = 57077
$ws.Range("D309").Value = 93.08
$ws.Range("E309").Value = 111.2
$ws.Range("F309").Value = 1
$ws.Range("G309").Value = 93.08

# Row 310
$ws.Range("B310").Value = 61610
$ws.Range("D310").Value = 102.71
$ws.Range("E310").Value = 122.71
$ws.Range("F310").Value = -58
$ws.Range("G310").Value = -5957.18

# Row 338
$ws.Range("B338").Value = 63520
$ws.Range("E338").Value = 153.4
$ws.Range("F338").Value = 97
$ws.Range("G338").Value = 13995.16

# Row 339
$ws.Range("B339").Value = 55373
$ws.Range("E339").Value = 163.62
$ws.Range("F339").Value = -94
$ws.Range("G339").Value = -13562.32

# Row 342
$ws.Range("B342").Value = 57802
$ws.Range("E342").Value = 162.71
$ws.Range("F342").Value = -79
$ws.Range("G342").Value = -11334.92

# Row 343
$ws.Range("B343").Value = 63571
$ws.Range("F343").Value = 29
$ws.Range("G343").Value = 4160.92

# Row 344
$ws.Range("B344").Value = 63531
$ws.Range("E344").Value = 152.53
$ws.Range("F344").Value = 80
$ws.Range("G344").Value = 11478.4

# Row 347
$ws.Range("B347").Value = 63510
$ws.Range("E347").Value = 50.66
$ws.Range("F347").Value = 167
$ws.Range("G347").Value = 7955.88

# Row 348
$ws.Range("B348").Value = 55356
$ws.Range("E348").Value = 54.04
$ws.Range("F348").Value = -158
$ws.Range("G348").Value = -7527.12

# Row 374
$ws.Range("B374").Value = 60325
$ws.Range("E374").Value = 151.57
$ws.Range("F374").Value = -102
$ws.Range("G374").Value = -12939.72

# Row 375
$ws.Range("B375").Value = 63560
$ws.Range("E375").Value = 134.87
$ws.Range("F375").Value = 104
$ws.Range("G375").Value = 13193.44

# Row 381
$ws.Range("B381").Value = 57817
$ws.Range("F381").Value = 3
$ws.Range("G381").Value = 239.43

# Row 382
$ws.Range("B382").Value = 62865
$ws.Range("F382").Value = 151
$ws.Range("G382").Value = 12051.31

# Row 392
$ws.Range("B392").Value = 57835
$ws.Range("F392").Value = 1
$ws.Range("G392").Value = 59.13

# Row 393
$ws.Range("B393").Value = 62933
$ws.Range("F393").Value = 146
$ws.Range("G393").Value = 8632.98

# Row 413
$ws.Range("B413").Value = 57857
$ws.Range("F413").Value = 3
$ws.Range("G413").Value = 453.51

# Row 414
$ws.Range("B414").Value = 63008
$ws.Range("F414").Value = 504
$ws.Range("G414").Value = 76189.67999999999

# Row 575
$ws.Range("B575").Value = 53263
$ws.Range("E575").Value = 15.29
$ws.Range("F575").Value = -309
$ws.Range("G575").Value = -3958.29

# Row 576
$ws.Range("B576").Value = 65066
$ws.Range("E576").Value = 13.61
$ws.Range("F576").Value = 313
$ws.Range("G576").Value = 4009.53

# Row 578
$ws.Range("B578").Value = 45695
$ws.Range("E578").Value = 23.58
$ws.Range("F578").Value = -36
$ws.Range("G578").Value = -710.28

# Row 579
$ws.Range("B579").Value = 64915
$ws.Range("E579").Value = 20.98
$ws.Range("F579").Value = 40
$ws.Range("G579").Value = 789.2

# Row 582
$ws.Range("B582").Value = 45706
$ws.Range("E582").Value = 23.58
$ws.Range("F582").Value = -202
$ws.Range("G582").Value = -3985.46

# Row 583
$ws.Range("B583").Value = 64922
$ws.Range("E583").Value = 20.98
$ws.Range("F583").Value = 207
$ws.Range("G583").Value = 4084.11

# Row 585
$ws.Range("B585").Value = 64927
$ws.Range("E585").Value = 17.26
$ws.Range("F585").Value = 295
$ws.Range("G585").Value = 4784.9

# Row 586
$ws.Range("B586").Value = 45718
$ws.Range("E586").Value = 19.38
$ws.Range("F586").Value = -294
$ws.Range("G586").Value = -4768.68

# Row 593
$ws.Range("B593").Value = 45702
$ws.Range("E593").Value = 31.43
$ws.Range("F593").Value = -215
$ws.Range("G593").Value = -5654.5

# Row 594
$ws.Range("B594").Value = 64919
$ws.Range("E594").Value = 27.97
$ws.Range("F594").Value = 224
$ws.Range("G594").Value = 5891.2

# Row 596
$ws.Range("B596").Value = 65067
$ws.Range("E596").Value = 15.65
$ws.Range("F596").Value = 338
$ws.Range("G596").Value = 4978.74

# Row 597
$ws.Range("B597").Value = 53595
$ws.Range("E597").Value = 17.61
$ws.Range("F597").Value = -335
$ws.Range("G597").Value = -4934.55

# Row 679
$ws.Range("B679").Value = 64810
$ws.Range("E679").Value = 291.22
$ws.Range("F679").Value = 7
$ws.Range("G679").Value = 1917.44

# Row 680
$ws.Range("B680").Value = 53319
$ws.Range("E680").Value = 310.64
$ws.Range("F680").Value = -6
$ws.Range("G680").Value = -1643.52

# Row 707
$ws.Range("B707").Value = 64836
$ws.Range("E707").Value = 104.71
$ws.Range("F707").Value = 7
$ws.Range("G707").Value = 689.5

# Row 708
$ws.Range("B708").Value = 60031
$ws.Range("E708").Value = 111.69
$ws.Range("F708").Value = -5
$ws.Range("G708").Value = -492.5

# Row 712
$ws.Range("B712").Value = 60022
$ws.Range("E712").Value = 37.22
$ws.Range("F712").Value = -113
$ws.Range("G712").Value = -3709.79

# Row 713
$ws.Range("B713").Value = 64830
$ws.Range("E713").Value = 34.9
$ws.Range("F713").Value = 117
$ws.Range("G713").Value = 3841.11

# Row 864
$ws.Range("B864").Value = 54751
$ws.Range("E864").Value = 46.34
$ws.Range("F864").Value = -19
$ws.Range("G864").Value = -776.53

# Row 865
$ws.Range("B865").Value = 65079
$ws.Range("E865").Value = 43.44
$ws.Range("F865").Value = 21
$ws.Range("G865").Value = 858.27

Write-Host "Applied 188 cell updates across 47 row groups."
